# Generate Report for Handoff
# Moves the localization status from "In Translation" to "Ready for
# handoff" and refreshes the associated handoff timestamps, mirroring a
# fresh run of the status-report generator.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: per-locale status + latest handoff-xliff timestamp ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-09-01 01:03:11"

# --- zh-cn detail sheet: status + latest handoff datetime ---
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-09-01 01:03:04"

# --- de-de detail sheet: status + latest handoff datetime ---
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-09-01 01:03:11"

# --- widen the Status columns to fit the longer "Ready for handoff" text ---
$wsOverview.Range("E1").ColumnWidth = 17.0
$wsOverview.Range("F1").ColumnWidth = 17.0
$wsZhCn.Range("C1").ColumnWidth = 17.0
$wsDeDe.Range("C1").ColumnWidth = 17.0
